# chore: update Sheets via scheduled runner
# Refreshes cached Market Board pricing (currentAveragePrice*) and the
# derived Leve price/profit columns (H:N) for affected leve rows across
# each job sheet of the workbook.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# row 4
$ws.Range("H4").Value = 606.7692
$ws.Range("I4").Value = 601.04
$ws.Range("K4").Value = 601.04
$ws.Range("M4").Value = -487.04
# row 43
$ws.Range("H43").Value = 6666.3335
$ws.Range("J43").Value = 6666.3335
$ws.Range("L43").Value = 6666.3335
$ws.Range("N43").Value = -6804.3335
# row 53
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = 0
# row 64
$ws.Range("H64").Value = 4251.6665
$ws.Range("I64").Value = 4123.75
$ws.Range("J64").Value = 5275
$ws.Range("K64").Value = 4123.75
$ws.Range("L64").Value = 5275
$ws.Range("M64").Value = -3875.75
$ws.Range("N64").Value = -5771
# row 67
$ws.Range("H67").Value = 4251.6665
$ws.Range("I67").Value = 4123.75
$ws.Range("J67").Value = 5275
$ws.Range("K67").Value = 4123.75
$ws.Range("L67").Value = 5275
$ws.Range("M67").Value = -3265.75
$ws.Range("N67").Value = -6991
# row 137
$ws.Range("H137").Value = 9169.700000000001
$ws.Range("I137").Value = 6099.5713
$ws.Range("K137").Value = 18298.7139
$ws.Range("M137").Value = -15748.7139

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 742.25
$ws.Range("I32").Value = 742.25
$ws.Range("K32").Value = 742.25
$ws.Range("M32").Value = -455.25
# row 61
$ws.Range("H61").Value = 5057.2856
$ws.Range("I61").Value = 5138.615
$ws.Range("K61").Value = 5138.615
$ws.Range("M61").Value = -4926.615
# row 63
$ws.Range("H63").Value = 6134.2856
$ws.Range("I63").Value = 4735
$ws.Range("K63").Value = 4735
$ws.Range("M63").Value = -4049
# row 66
$ws.Range("H66").Value = 6134.2856
$ws.Range("I66").Value = 4735
$ws.Range("K66").Value = 23675
$ws.Range("M66").Value = -20243
# row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# row 136
$ws.Range("H136").Value = 5057.2856
$ws.Range("I136").Value = 5138.615
$ws.Range("K136").Value = 15415.845
$ws.Range("M136").Value = -12865.845

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 699.6667
$ws.Range("J86").Value = 549
$ws.Range("L86").Value = 549
$ws.Range("N86").Value = -2795
# row 89
$ws.Range("H89").Value = 699.6667
$ws.Range("J89").Value = 549
$ws.Range("L89").Value = 2745
$ws.Range("N89").Value = -13977
# row 61
$ws.Range("H61").Value = 1205
$ws.Range("I61").Value = 1205
$ws.Range("K61").Value = 1205
$ws.Range("M61").Value = -1003
# row 113
$ws.Range("H113").Value = 1205
$ws.Range("I113").Value = 1205
$ws.Range("K113").Value = 1205
$ws.Range("M113").Value = 965

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value = 397.16666
$ws.Range("I2").Value = 64.333336
$ws.Range("J2").Value = 730
$ws.Range("K2").Value = 64.333336
$ws.Range("L2").Value = 730
$ws.Range("M2").Value = 48.666664
$ws.Range("N2").Value = -956
# row 15
$ws.Range("H15").Value = 985
$ws.Range("J15").Value = 798.75
$ws.Range("L15").Value = 798.75
$ws.Range("N15").Value = -1138.75
# row 37
$ws.Range("H37").Value = 15000
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
# row 60
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 50000
$ws.Range("K60").Value = 50000
$ws.Range("M60").Value = -49489
# row 132
$ws.Range("H132").Value = 4374.875
$ws.Range("I132").Value = 999.8
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 2999.4
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -469.3999999999996
$ws.Range("N132").Value = -35060
# row 134
$ws.Range("H134").Value = 3364.125
$ws.Range("I134").Value = 984.3333
$ws.Range("K134").Value = 2952.9999
$ws.Range("M134").Value = -417.9998999999998

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# row 44
$ws.Range("H44").Value = 499
$ws.Range("I44").Value = 498
$ws.Range("J44").Value = 499.5
$ws.Range("K44").Value = 1494
$ws.Range("L44").Value = 1498.5
$ws.Range("M44").Value = -1096
$ws.Range("N44").Value = -2294.5
# row 114
$ws.Range("H114").Value = 1120.6666
$ws.Range("I114").Value = 300
$ws.Range("J114").Value = 1531
$ws.Range("K114").Value = 900
$ws.Range("L114").Value = 4593
$ws.Range("M114").Value = 2354
$ws.Range("N114").Value = -11101
# row 132
$ws.Range("H132").Value = 760
$ws.Range("J132").Value = 395
$ws.Range("L132").Value = 3555
$ws.Range("N132").Value = -8615

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3250
$ws.Range("I102").Value = 2500
$ws.Range("K102").Value = 2500
$ws.Range("M102").Value = -878
# row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
# row 126
$ws.Range("H126").Value = 7757
$ws.Range("J126").Value = 7757
$ws.Range("L126").Value = 23271
$ws.Range("N126").Value = -28211
# row 132
$ws.Range("H132").Value = 4981.7646
$ws.Range("I132").Value = 2833.25
$ws.Range("K132").Value = 8499.75
$ws.Range("M132").Value = -5969.75

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# row 42
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("N42").Value = 0
# row 49
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
# row 55
$ws.Range("H55").Value = 1029.8
$ws.Range("I55").Value = 912.25
$ws.Range("K55").Value = 912.25
$ws.Range("M55").Value = -739.25

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 2000.8572
$ws.Range("I122").Value = 2000.8572
$ws.Range("K122").Value = 6002.571599999999
$ws.Range("M122").Value = -3552.571599999999
# row 132
$ws.Range("H132").Value = 11564.167
$ws.Range("I132").Value = 10396.143
$ws.Range("K132").Value = 31188.429
$ws.Range("M132").Value = -28658.429
